$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 434.48148
$ws.Range("I15").Value = 434.48148
$ws.Range("K15").Value = 1303.44444
$ws.Range("M15").Value = -1134.44444
$ws.Range("H40").Value = 1879.9
$ws.Range("I40").Value = 1849.875
$ws.Range("K40").Value = 1849.875
$ws.Range("M40").Value = -1674.875
$ws.Range("H43").Value = 5975
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5138
$ws.Range("H70").Value = 119494.836
$ws.Range("I70").Value = 3250
$ws.Range("J70").Value = 142743.8
$ws.Range("K70").Value = 9750
$ws.Range("L70").Value = 428231.4
$ws.Range("M70").Value = -9480
$ws.Range("N70").Value = -428771.4
$ws.Range("H73").Value = 119494.836
$ws.Range("I73").Value = 3250
$ws.Range("J73").Value = 142743.8
$ws.Range("K73").Value = 9750
$ws.Range("L73").Value = 428231.4
$ws.Range("M73").Value = -8814
$ws.Range("N73").Value = -430103.4
$ws.Range("H100").Value = 1319.6522
$ws.Range("I100").Value = 947.3333
$ws.Range("K100").Value = 947.3333
$ws.Range("M100").Value = -406.3333
$ws.Range("H138").Value = 12096.395
$ws.Range("J138").Value = 13548.577
$ws.Range("L138").Value = 40645.731
$ws.Range("N138").Value = -50925.731

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H32").Value = 31319
$ws.Range("I32").Value = 33284.332
$ws.Range("K32").Value = 33284.332
$ws.Range("M32").Value = -32997.332
$ws.Range("H37").Value = 13333
$ws.Range("J37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("N37").Value = -10546
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 1650
$ws.Range("I61").Value = 975
$ws.Range("K61").Value = 975
$ws.Range("M61").Value = -763
$ws.Range("H80").Value = 89998.5
$ws.Range("J80").Value = 89998.5
$ws.Range("L80").Value = 89998.5
$ws.Range("N80").Value = -91994.5
$ws.Range("H83").Value = 89998.5
$ws.Range("J83").Value = 89998.5
$ws.Range("L83").Value = 269995.5
$ws.Range("N83").Value = -279979.5
$ws.Range("H88").Value = 3361.4
$ws.Range("J88").Value = 3451.75
$ws.Range("L88").Value = 3451.75
$ws.Range("N88").Value = -4263.75
$ws.Range("H91").Value = 3361.4
$ws.Range("J91").Value = 3451.75
$ws.Range("L91").Value = 3451.75
$ws.Range("N91").Value = -6259.75
$ws.Range("H132").Value = 1585.7222
$ws.Range("I132").Value = 1416.2667
$ws.Range("K132").Value = 4248.800099999999
$ws.Range("M132").Value = -1718.800099999999
$ws.Range("H136").Value = 1650
$ws.Range("I136").Value = 975
$ws.Range("K136").Value = 2925
$ws.Range("M136").Value = -375

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2336.926
$ws.Range("I20").Value = 1918.8422
$ws.Range("K20").Value = 1918.8422
$ws.Range("M20").Value = -1671.8422
$ws.Range("H34").Value = 8000
$ws.Range("I34").Value = 8000
$ws.Range("K34").Value = 8000
$ws.Range("M34").Value = -7886
$ws.Range("H86").Value = 7519.4443
$ws.Range("I86").Value = 7499.75
$ws.Range("J86").Value = 7535.2
$ws.Range("K86").Value = 7499.75
$ws.Range("L86").Value = 7535.2
$ws.Range("M86").Value = -6376.75
$ws.Range("N86").Value = -9781.200000000001
$ws.Range("H89").Value = 7519.4443
$ws.Range("I89").Value = 7499.75
$ws.Range("J89").Value = 7535.2
$ws.Range("K89").Value = 37498.75
$ws.Range("L89").Value = 37676
$ws.Range("M89").Value = -31882.75
$ws.Range("N89").Value = -48908
$ws.Range("H105").Value = 4601.826
$ws.Range("I105").Value = 4087.923
$ws.Range("J105").Value = 5269.9
$ws.Range("K105").Value = 4087.923
$ws.Range("L105").Value = 5269.9
$ws.Range("M105").Value = -2340.923
$ws.Range("N105").Value = -8763.9

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13874.333
$ws.Range("I99").Value = 9664.571
$ws.Range("J99").Value = 16553.273
$ws.Range("K99").Value = 9664.571
$ws.Range("L99").Value = 16553.273
$ws.Range("M99").Value = -8166.571
$ws.Range("N99").Value = -19549.273
$ws.Range("H105").Value = 6497.1665
$ws.Range("I105").Value = 4661.6665
$ws.Range("J105").Value = 8332.666999999999
$ws.Range("K105").Value = 4661.6665
$ws.Range("L105").Value = 8332.666999999999
$ws.Range("M105").Value = -2914.6665
$ws.Range("N105").Value = -11826.667
$ws.Range("H122").Value = 1734.75
$ws.Range("J122").Value = 1399.8
$ws.Range("L122").Value = 4199.4
$ws.Range("N122").Value = -9099.4
$ws.Range("H126").Value = 13874.333
$ws.Range("I126").Value = 9664.571
$ws.Range("J126").Value = 16553.273
$ws.Range("K126").Value = 28993.713
$ws.Range("L126").Value = 49659.819
$ws.Range("M126").Value = -26523.713
$ws.Range("N126").Value = -54599.819
$ws.Range("H134").Value = 3787.9285
$ws.Range("I134").Value = 3166.111
$ws.Range("K134").Value = 9498.332999999999
$ws.Range("M134").Value = -6963.332999999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4299.6
$ws.Range("I80").Value = 4499.6665
$ws.Range("J80").Value = 3999.5
$ws.Range("K80").Value = 13498.9995
$ws.Range("L80").Value = 11998.5
$ws.Range("M80").Value = -12562.9995
$ws.Range("N80").Value = -13870.5
$ws.Range("H83").Value = 4299.6
$ws.Range("I83").Value = 4499.6665
$ws.Range("J83").Value = 3999.5
$ws.Range("K83").Value = 40496.9985
$ws.Range("L83").Value = 35995.5
$ws.Range("M83").Value = -35816.9985
$ws.Range("N83").Value = -45355.5
$ws.Range("H125").Value = 500
$ws.Range("I125").Value = 500
$ws.Range("K125").Value = 1500
$ws.Range("M125").Value = 3420
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 12000
$ws.Range("J132").Value = 12000
$ws.Range("L132").Value = 108000
$ws.Range("N132").Value = -113060

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 268
$ws.Range("I3").Value = 300
$ws.Range("J3").Value = 252
$ws.Range("K3").Value = 300
$ws.Range("L3").Value = 252
$ws.Range("M3").Value = -184
$ws.Range("N3").Value = -484
$ws.Range("H57").Value = 42499.5
$ws.Range("J57").Value = 59999
$ws.Range("L57").Value = 59999
$ws.Range("N57").Value = -61639
$ws.Range("H70").Value = 1949.5
$ws.Range("I70").Value = 1949.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 1949.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -1679.5
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 1949.5
$ws.Range("I73").Value = 1949.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 1949.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1013.5
$ws.Range("N73").ClearContents()
$ws.Range("H97").Value = 2533.6
$ws.Range("I97").Value = 2535.2856
$ws.Range("K97").Value = 2535.2856
$ws.Range("M97").Value = -2039.2856
$ws.Range("H122").Value = 383401.2
$ws.Range("I122").Value = 70459.92999999999
$ws.Range("J122").Value = 718695.4
$ws.Range("K122").Value = 211379.79
$ws.Range("L122").Value = 2156086.2
$ws.Range("M122").Value = -208929.79
$ws.Range("N122").Value = -2160986.2
$ws.Range("H132").Value = 4749.6665
$ws.Range("I132").Value = 3972.2727
$ws.Range("K132").Value = 11916.8181
$ws.Range("M132").Value = -9386.8181

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1833.2222
$ws.Range("J22").Value = 2266.6667
$ws.Range("L22").Value = 2266.6667
$ws.Range("N22").Value = -2856.6667
$ws.Range("H27").Value = 1833.2222
$ws.Range("J27").Value = 2266.6667
$ws.Range("L27").Value = 2266.6667
$ws.Range("N27").Value = -2480.6667
$ws.Range("H132").Value = 4752
$ws.Range("I132").Value = 4096.5
$ws.Range("K132").Value = 12289.5
$ws.Range("M132").Value = -9759.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 6781.25
$ws.Range("J40").Value = 6781.25
$ws.Range("L40").Value = 6781.25
$ws.Range("N40").Value = -7079.25
$ws.Range("H122").Value = 1502
$ws.Range("I122").Value = 1502
$ws.Range("K122").Value = 4506
$ws.Range("M122").Value = -2056
